$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values stay stored as text, matching the source data
# (some prices like "1.005" would otherwise be auto-converted to numbers by Excel)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.968.25"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.827.55"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "311.04"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").Value = "0.4625"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").Value = "0.3686"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").Value = "0.07329"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").Value = "0.8763"
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("D11").Value = "0.07927"
$ws.Range("E11").Value = "  +5.19%  "
$ws.Range("D12").Value = "19.62"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "1.807.90"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "5.328"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").Value = "6.540"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "91.16"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "0.000008921"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "14.72"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").Value = "26.978.50"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").Value = "10.53"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "2.023.10"
$ws.Range("E24").Value = "  -3.71%  "
$ws.Range("D25").Value = "152.40"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "1.848"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").Value = "2.024"
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("D29").Value = "5.081"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "115.59"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").Value = "0.08847"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "2.966"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "0.7326"
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("D34").Value = "4.433"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("D35").Value = "1.129"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "1.073"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "2.442"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "0.01940"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "0.05233"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").Value = "2.954"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("D41").Value = "6.963"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").Value = "0.5125"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "0.1622"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "8.091"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.4792"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").Value = "1.004"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "10.17"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").Value = "101.64"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").Value = "1.614"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D50").Value = "0.06199"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "64.36"
$ws.Range("E51").Value = "  +0.81%  "

# Restore the default cell style so no stray formatting is introduced
$ws.Range("D2:D51").Style = "Normal"
